$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remember old values before overwriting anything
$oldO1 = $ws.Range("O1").Value()
$oldO2 = $ws.Range("O2").Value()
$oldO2NumberFormat = $ws.Range("O2").NumberFormat()
$oldO2VerticalAlignment = $ws.Range("O2").VerticalAlignment()

# Row 1 (header) edits
$ws.Range("O1").Value = "TEXT4"
$ws.Range("P1").Value = "TEXT5"
$ws.Range("Q1").Value = $oldO1

# Row 2 (data) edits: move "ReportProductHolding.xlsx" from O2 to Q2
$ws.Range("Q2").Value = $oldO2
$ws.Range("Q2").NumberFormat = $oldO2NumberFormat
$ws.Range("Q2").VerticalAlignment = $oldO2VerticalAlignment

$ws.Range("O2").ClearContents()
$ws.Range("O2").Style = $ws.Range("A1").Style()

# Column widths (O & P take the narrow "TEXT" width, Q takes the wide width
# previously used by O). The underlying engine quantizes ColumnWidth to
# whole-pixel steps, so these inputs are chosen to land as close as possible
# to the exact target stored widths (6.140625 and 25.5703125).
$ws.Columns.Item(15).ColumnWidth = 5.333333333333333
$ws.Columns.Item(16).ColumnWidth = 5.333333333333333
$ws.Columns.Item(17).ColumnWidth = 24.666666666666668

# Update selection to Q2
$ws.Range("Q2").Select()
